# Refresh Universalis-derived market price & profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across all Leve sheets.
# Source data pulled from the latest Universalis snapshot by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 2000
$ws.Range("J44").Value = 2000
$ws.Range("L44").Value = 2000
$ws.Range("N44").Value = -2924
$ws.Range("H125").Value = 2000
$ws.Range("J125").Value = 2000
$ws.Range("L125").Value = 18000
$ws.Range("N125").Value = -22920
$ws.Range("H127").Value = 3249.75
$ws.Range("J127").Value = 3249.75
$ws.Range("L127").Value = 9749.25
$ws.Range("N127").Value = -19669.25
$ws.Range("H132").Value = 1170.0769
$ws.Range("I132").Value = 1170.0769
$ws.Range("K132").Value = 3510.2307
$ws.Range("M132").Value = -980.2307000000001
$ws.Range("H138").Value = 3276.1091
$ws.Range("I138").Value = 3592.6
$ws.Range("K138").Value = 10777.8
$ws.Range("M138").Value = -5637.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1241.2354
$ws.Range("J2").Value = 1452.2
$ws.Range("L2").Value = 1452.2
$ws.Range("N2").Value = -1678.2
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H32").Value = 7176.108
$ws.Range("I32").Value = 4694.0967
$ws.Range("K32").Value = 4694.0967
$ws.Range("M32").Value = -4407.0967
$ws.Range("H61").Value = 2582.125
$ws.Range("I61").Value = 1683.5
$ws.Range("J61").Value = 2881.6667
$ws.Range("K61").Value = 1683.5
$ws.Range("L61").Value = 2881.6667
$ws.Range("M61").Value = -1471.5
$ws.Range("N61").Value = -3305.6667
$ws.Range("H102").Value = 2293.5
$ws.Range("I102").Value = 2338.3635
$ws.Range("J102").Value = 1800
$ws.Range("K102").Value = 2338.3635
$ws.Range("L102").Value = 1800
$ws.Range("M102").Value = -716.3634999999999
$ws.Range("N102").Value = -5044
$ws.Range("H116").Value = 1241.2354
$ws.Range("J116").Value = 1452.2
$ws.Range("L116").Value = 1452.2
$ws.Range("N116").Value = -6040.2
$ws.Range("H120").Value = 29710
$ws.Range("J120").Value = 29710
$ws.Range("L120").Value = 29710
$ws.Range("N120").Value = -39386
$ws.Range("H122").Value = 2899.4
$ws.Range("I122").Value = 2624.25
$ws.Range("K122").Value = 7872.75
$ws.Range("M122").Value = -5422.75
$ws.Range("H136").Value = 2582.125
$ws.Range("I136").Value = 1683.5
$ws.Range("J136").Value = 2881.6667
$ws.Range("K136").Value = 5050.5
$ws.Range("L136").Value = 8645.000100000001
$ws.Range("M136").Value = -2500.5
$ws.Range("N136").Value = -13745.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1241.2354
$ws.Range("J3").Value = 1452.2
$ws.Range("L3").Value = 1452.2
$ws.Range("N3").Value = -1680.2
$ws.Range("H107").Value = 625.8889
$ws.Range("I107").Value = 607.5714
$ws.Range("K107").Value = 607.5714
$ws.Range("M107").Value = 1312.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 499
$ws.Range("I7").Value = 516.8
$ws.Range("J7").Value = 469.33334
$ws.Range("K7").Value = 516.8
$ws.Range("L7").Value = 469.33334
$ws.Range("M7").Value = -403.8
$ws.Range("N7").Value = -695.33334
$ws.Range("H31").Value = 1599
$ws.Range("I31").Value = 1599.5
$ws.Range("J31").Value = 1598
$ws.Range("K31").Value = 1599.5
$ws.Range("L31").Value = 1598
$ws.Range("M31").Value = -1304.5
$ws.Range("N31").Value = -2188
$ws.Range("H34").Value = 1599
$ws.Range("I34").Value = 1599.5
$ws.Range("J34").Value = 1598
$ws.Range("K34").Value = 1599.5
$ws.Range("L34").Value = 1598
$ws.Range("M34").Value = -1397.5
$ws.Range("N34").Value = -2002
$ws.Range("H86").Value = 20477.75
$ws.Range("I86").Value = 9857
$ws.Range("J86").Value = 29167.455
$ws.Range("K86").Value = 9857
$ws.Range("L86").Value = 29167.455
$ws.Range("M86").Value = -8734
$ws.Range("N86").Value = -31413.455
$ws.Range("H89").Value = 20477.75
$ws.Range("I89").Value = 9857
$ws.Range("J89").Value = 29167.455
$ws.Range("K89").Value = 49285
$ws.Range("L89").Value = 145837.275
$ws.Range("M89").Value = -43669
$ws.Range("N89").Value = -157069.275
$ws.Range("H105").Value = 1764.8334
$ws.Range("I105").Value = 1857.8
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 1857.8
$ws.Range("L105").Value = 1300
$ws.Range("M105").Value = -110.8
$ws.Range("N105").Value = -4794
$ws.Range("H107").Value = 1934.4117
$ws.Range("I107").Value = 1030
$ws.Range("K107").Value = 1030
$ws.Range("M107").Value = 890

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 170832.33
$ws.Range("J80").Value = 203999.6
$ws.Range("L80").Value = 611998.8
$ws.Range("N80").Value = -613870.8
$ws.Range("H83").Value = 170832.33
$ws.Range("J83").Value = 203999.6
$ws.Range("L83").Value = 1835996.4
$ws.Range("N83").Value = -1845356.4
$ws.Range("H125").Value = 7499
$ws.Range("J125").Value = 7499
$ws.Range("L125").Value = 22497
$ws.Range("N125").Value = -32337

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5400.2
$ws.Range("I80").Value = 4248.75
$ws.Range("J80").Value = 10006
$ws.Range("K80").Value = 4248.75
$ws.Range("L80").Value = 10006
$ws.Range("M80").Value = -3250.75
$ws.Range("N80").Value = -12002
$ws.Range("H83").Value = 5400.2
$ws.Range("I83").Value = 4248.75
$ws.Range("J83").Value = 10006
$ws.Range("K83").Value = 21243.75
$ws.Range("L83").Value = 50030
$ws.Range("M83").Value = -16251.75
$ws.Range("N83").Value = -60014
$ws.Range("H113").Value = 1413.3
$ws.Range("I113").Value = 1413.3
$ws.Range("K113").Value = 1413.3
$ws.Range("M113").Value = 756.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 19180
$ws.Range("J81").Value = 19180
$ws.Range("L81").Value = 19180
$ws.Range("N81").Value = -21176
$ws.Range("H84").Value = 19180
$ws.Range("J84").Value = 19180
$ws.Range("L84").Value = 57540
$ws.Range("N84").Value = -67524
$ws.Range("H93").Value = 3158.4
$ws.Range("I93").Value = 3698.75
$ws.Range("J93").Value = 997
$ws.Range("K93").Value = 3698.75
$ws.Range("L93").Value = 997
$ws.Range("M93").Value = -2450.75
$ws.Range("N93").Value = -3493
$ws.Range("H132").Value = 3705.125
$ws.Range("I132").Value = 2940.5
$ws.Range("K132").Value = 8821.5
$ws.Range("M132").Value = -6291.5
$ws.Range("H136").Value = 4343.5
$ws.Range("I136").Value = 4678.2856
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 14034.8568
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -11484.8568
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10046.25
$ws.Range("I81").Value = 4677
$ws.Range("K81").Value = 9354
$ws.Range("M81").Value = -8293
$ws.Range("H84").Value = 10046.25
$ws.Range("I84").Value = 4677
$ws.Range("K84").Value = 46770
$ws.Range("M84").Value = -41466
$ws.Range("H113").Value = 186.4
$ws.Range("I113").Value = 143.5
$ws.Range("K113").Value = 430.5
$ws.Range("M113").Value = 1739.5
$ws.Range("H132").Value = 1770.4348
$ws.Range("I132").Value = 1581.2354
$ws.Range("K132").Value = 4743.706200000001
$ws.Range("M132").Value = -2213.706200000001
$ws.Range("H136").Value = 2223.7576
$ws.Range("J136").Value = 3289.9
$ws.Range("L136").Value = 9869.700000000001
$ws.Range("N136").Value = -14969.7
